# Generate Report for Handback
#
# This script reorders the two "in sync with en-US" handback rows
# (fffffc6ff801-...md and ffffffde59b5f6-...md) so that the row for
# ffffffde59b5f6-...md comes first, flips the 52595375-...md row's
# status from "Ready for handoff" to "Handed back: in sync with en-US",
# and refreshes the handback timestamps for that row on the zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAt($worksheet, $addr) {
    foreach ($hl in $worksheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl
        }
    }
    return $null
}

function Swap-HyperlinkCell($worksheet, $addr1, $addr2) {
    # Swap the displayed text (and underlying cell text) between two
    # single-cell hyperlinked ranges, while leaving each hyperlink's
    # underlying target relationship (Address/rId) untouched - this is
    # what actually happened in the source edit (only the display text
    # moved between rows, the relationship ids were not renumbered).
    $hl1 = Get-HyperlinkAt $worksheet $addr1
    $hl2 = Get-HyperlinkAt $worksheet $addr2

    $textValue1 = $hl1.TextToDisplay
    $textValue2 = $hl2.TextToDisplay

    $hl1.TextToDisplay = $textValue2
    $worksheet.Range($addr1).Value = $textValue2

    $hl2.TextToDisplay = $textValue1
    $worksheet.Range($addr2).Value = $textValue1
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
Swap-HyperlinkCell $wsOverview "`$A`$2" "`$A`$3"
$wsOverview.Range("B4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C4").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Swap-HyperlinkCell $wsZhCn "`$A`$2" "`$A`$3"
$wsZhCn.Range("B4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G4").Value = "2016-03-03 15:44:02"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
Swap-HyperlinkCell $wsDeDe "`$A`$2" "`$A`$3"
$wsDeDe.Range("B4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G4").Value = "2016-03-03 15:44:25"
